# Generate Report for Handback
#
# - Flips the "Status" column (shared by the Overview summary sheet and the
#   per-locale detail sheets) from "Ready for handoff" to
#   "Handed back: in sync with en-US".
# - Populates the "Latest Target File" / "Latest Handback File" columns
#   (F/G) for the two data rows on each locale sheet, with hyperlinks that
#   mirror the existing Source/Handoff-file hyperlinks.
# - Records the handback timestamp in the "Latest Handback DateTime"
#   column (H) per locale.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAddress($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl.Address
        }
    }
    return $null
}

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status shown for both locales, both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn detail sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Reuse the targets already used by the Source File Name (A) and
# Latest Handoff File (D) hyperlinks for the new Latest Target File (F)
# and Latest Handback File (G) columns.
$zhMdAddress = Get-HyperlinkAddress $wsZh "`$A`$2"
$zhXlfAddress = Get-HyperlinkAddress $wsZh "`$D`$2"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhMdAddress, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfAddress, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhMdAddress, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfAddress, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# Latest Handback DateTime
$wsZh.Range("H2").Value = "2016-03-23 06:35:01"
$wsZh.Range("H3").Value = "2016-03-23 06:35:01"

# --- de-de detail sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$deMdAddress = Get-HyperlinkAddress $wsDe "`$A`$2"
$deXlfAddress = Get-HyperlinkAddress $wsDe "`$D`$2"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deMdAddress, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfAddress, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deMdAddress, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfAddress, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")

# Latest Handback DateTime
$wsDe.Range("H2").Value = "2016-03-23 06:35:12"
$wsDe.Range("H3").Value = "2016-03-23 06:35:12"

Write-Output "Handback report generated"
